# feat: add 2022-Q1 data
#
# - Insert a new worksheet "2022-Q1" (holding positions) right before the
#   "总计" (grand-total) sheet, matching the layout used by "2021-Q2" /
#   "2021-Q3".
# - Add a new leading row to the "总计" sheet summarising the new quarter.

$wb = $excel.ActiveWorkbook

$q3 = $wb.Worksheets.Item("2021-Q3")
$totalBefore = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------
# 1. Create the "2022-Q1" worksheet, positioned right before "总计".
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Add($totalBefore)
$q1.Name = "2022-Q1"

# NOTE: inserting a sheet shifts the positional index of every sheet
# after it, and worksheet variables here are resolved positionally, so
# any reference obtained *before* the Add() call (e.g. $totalBefore) can
# silently point at the wrong sheet afterwards. Always re-fetch "总计"
# by name once the sheet collection has changed.
$total = $wb.Worksheets.Item("总计")

# Copy header row (values + formatting) from the "2021-Q3" sheet so the
# new sheet matches the existing per-quarter layout exactly.
$q3.Range("B1:H1").Copy($q1.Range("B1"))
$q3.Range("A2").Copy($q1.Range("A2"))

$q1.Range("A2").Value = 0
# Use a leading apostrophe so these numeric-looking values are stored as
# text (matching the source data), not auto-converted to numbers.
$q1.Range("B2").Value = "'968029"
$q1.Range("C2").Value = "恒生指数基金M类人民币（对冲）份额"
$q1.Range("D2").Value = "'25.09"
$q1.Range("E2").Value = "'97.94"
$q1.Range("F2").Value = "'4.93"
$q1.Range("G2").Value = "'1.2369"
$q1.Range("H2").Value = 6

# ---------------------------------------------------------------------
# 2. Insert a new leading data row in "总计" for the 2022-Q1 summary,
#    shifting the existing rows down.
# ---------------------------------------------------------------------
$total.Rows.Item(2).Insert(-4121)
$total.Range("A2:D2").ClearFormats()

$total.Range("A4").Copy($total.Range("A2"))
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 1
$total.Range("D2").Value = "1.24"

$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
